$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2005899705014749
$ws.Range("C2").Value = 0.5132743362831859
$ws.Range("J2").Value = 0.02359882005899705
$ws.Range("O2").Value = 0.001474926253687316
$ws.Range("P2").Value = 0.168141592920354
$ws.Range("S2").Value = 0.09292035398230089
$ws.Range("B3").Value = 0.01466992665036675
$ws.Range("C3").Value = 0.02689486552567237
$ws.Range("J3").Value = 0.05134474327628362
$ws.Range("P3").Value = 0.7408312958435208
$ws.Range("S3").Value = 0.1662591687041565
$ws.Range("J4").Value = 0.1142857142857143
$ws.Range("P4").Value = 0.580952380952381
$ws.Range("S4").Value = 0.3047619047619048
$ws.Range("B6").Value = 0.07276507276507277
$ws.Range("D6").Value = 0.01247401247401247
$ws.Range("E6").Value = 0.002079002079002079
$ws.Range("F6").Value = 0.06029106029106029
$ws.Range("J6").Value = 0.2765072765072765
$ws.Range("O6").Value = 0.02494802494802495
$ws.Range("Q6").Value = 0.1891891891891892
$ws.Range("R6").Value = 0.04365904365904366
$ws.Range("S6").Value = 0.3180873180873181
$ws.Range("B7").Value = 0.1341463414634146
$ws.Range("D7").Value = 0.007317073170731708
$ws.Range("F7").Value = 0.06829268292682927
$ws.Range("J7").Value = 0.1317073170731707
$ws.Range("O7").Value = 0.02682926829268293
$ws.Range("Q7").Value = 0.148780487804878
$ws.Range("R7").Value = 0.08048780487804878
$ws.Range("S7").Value = 0.4024390243902439
$ws.Range("B8").Value = 0.1023359288097887
$ws.Range("D8").Value = 0.02335928809788654
$ws.Range("E8").Value = 0.001112347052280311
$ws.Range("F8").Value = 0.06451612903225806
$ws.Range("J8").Value = 0.1112347052280311
$ws.Range("O8").Value = 0.02558398220244716
$ws.Range("Q8").Value = 0.203559510567297
$ws.Range("R8").Value = 0.07341490545050056
$ws.Range("S8").Value = 0.3948832035595106
$ws.Range("B9").Value = 0.08306709265175719
$ws.Range("D9").Value = 0.02236421725239617
$ws.Range("F9").Value = 0.06070287539936102
$ws.Range("J9").Value = 0.1246006389776358
$ws.Range("O9").Value = 0.03833865814696485
$ws.Range("Q9").Value = 0.1789137380191693
$ws.Range("R9").Value = 0.08626198083067092
$ws.Range("S9").Value = 0.4057507987220447
$ws.Range("B10").Value = 0.1102507374631268
$ws.Range("D10").Value = 0.02359882005899705
$ws.Range("E10").Value = 0.001843657817109145
$ws.Range("F10").Value = 0.06452802359882005
$ws.Range("J10").Value = 0.1338495575221239
$ws.Range("O10").Value = 0.02949852507374631
$ws.Range("Q10").Value = 0.2367256637168142
$ws.Range("R10").Value = 0.06895280235988201
$ws.Range("S10").Value = 0.3307522123893805
$ws.Range("G11").Value = 0.1390625
$ws.Range("J11").Value = 0.1015625
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.546875
$ws.Range("S11").Value = 0.0125
$ws.Range("G12").Value = 0.7092391304347826
$ws.Range("J12").Value = 0.2146739130434783
$ws.Range("K12").Value = 0.0108695652173913
$ws.Range("L12").Value = 0.02989130434782609
$ws.Range("S12").Value = 0.03532608695652174
$ws.Range("G13").Value = 0.7727272727272727
$ws.Range("J13").Value = 0.125
$ws.Range("S13").Value = 0.1022727272727273
$ws.Range("F15").Value = 0.03379721669980119
$ws.Range("H15").Value = 0.1590457256461233
$ws.Range("I15").Value = 0.04572564612326044
$ws.Range("J15").Value = 0.341948310139165
$ws.Range("K15").Value = 0.07554671968190854
$ws.Range("M15").Value = 0.009940357852882704
$ws.Range("N15").Value = 0.003976143141153081
$ws.Range("O15").Value = 0.05168986083499006
$ws.Range("S15").Value = 0.2783300198807157
$ws.Range("F16").Value = 0.01986754966887417
$ws.Range("H16").Value = 0.1721854304635762
$ws.Range("I16").Value = 0.04856512141280353
$ws.Range("J16").Value = 0.3973509933774834
$ws.Range("K16").Value = 0.1125827814569536
$ws.Range("M16").Value = 0.01324503311258278
$ws.Range("O16").Value = 0.05077262693156733
$ws.Range("S16").Value = 0.1854304635761589
$ws.Range("F17").Value = 0.02358490566037736
$ws.Range("H17").Value = 0.1792452830188679
$ws.Range("I17").Value = 0.0660377358490566
$ws.Range("J17").Value = 0.430188679245283
$ws.Range("K17").Value = 0.09811320754716982
$ws.Range("M17").Value = 0.0169811320754717
$ws.Range("N17").Value = 0.0009433962264150943
$ws.Range("O17").Value = 0.06509433962264151
$ws.Range("S17").Value = 0.119811320754717
$ws.Range("F18").Value = 0.01724137931034483
$ws.Range("H18").Value = 0.1379310344827586
$ws.Range("I18").Value = 0.08908045977011494
$ws.Range("J18").Value = 0.4310344827586207
$ws.Range("K18").Value = 0.1005747126436782
$ws.Range("M18").Value = 0.005747126436781609
$ws.Range("N18").Value = 0.002873563218390805
$ws.Range("O18").Value = 0.08908045977011494
$ws.Range("S18").Value = 0.1264367816091954
$ws.Range("F19").Value = 0.02463651050080775
$ws.Range("H19").Value = 0.2079967689822294
$ws.Range("I19").Value = 0.06744749596122779
$ws.Range("J19").Value = 0.3760096930533118
$ws.Range("K19").Value = 0.1215670436187399
$ws.Range("M19").Value = 0.02382875605815832
$ws.Range("N19").Value = 0.0008077544426494346
$ws.Range("O19").Value = 0.07269789983844911
$ws.Range("S19").Value = 0.1050080775444265
